$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume 1h) hold numeric-/percent-looking strings that
# must stay stored as literal text (matching the source inlineStr cells), so force
# the Text number format before assigning - otherwise Excel auto-converts them to
# a number/percentage, same as typing them into a live sheet.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "324.79"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-2.49%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "44.46"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0.48%"

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-6.08%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08047"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-3.52%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "8.642"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-2.01%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.905"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-3.76%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "4.282"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-4.86%"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.698"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-6.90%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9413"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.73%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1162"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-6.91%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1863"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-4.24%"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.1005"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "5.03%"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.04357"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "10.72%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.1064"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.25%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001278"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-2.00%"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.005923"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-2.67%"

$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.592"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "2.44%"

$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3492"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.51%"

$ws.Range("B20").Value = "MCDex"
$ws.Range("C20").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.547"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-4.72%"

$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1371"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-0.07%"

$ws.Range("B22").Value = "ZBToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2532"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-1.54%"

$ws.Range("B23").Value = "CoinExToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04241"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-3.91%"

$ws.Range("B24").Value = "BitKan"
$ws.Range("C24").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001236"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-1.81%"

$ws.Range("B25").Value = "HotbitToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004547"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "3.19%"

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.85%"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003993"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-0.02%"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02627"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-6.08%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05472"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-3.99%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007695"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-3.16%"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1393"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-2.28%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007037"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-21.78%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002058"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-2.18%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.008677"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-17.36%"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00007103"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-2.32%"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.08%"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003644"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "12.02%"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002272"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.35%"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002102"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.08%"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002002"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.08%"
